$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are preserved as text, matching the
# original inline-string cell type (avoids Excel auto-converting
# numeric-looking strings like "8.870" into numbers and dropping
# formatting such as trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.207.69"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.852.75"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "0.6987"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").Value = "237.61"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.07867"
$ws.Range("E8").Value = "  +0.74%  "
$ws.Range("D9").Value = "0.3016"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").Value = "23.81"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").Value = "0.08131"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "1.849.51"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "5.184"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "0.7055"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "89.46"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "29.209.17"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "5.802"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "0.000007831"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").Value = "13.21"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "235.95"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "2.097.20"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").Value = "0.9997"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "7.509"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").Value = "162.54"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").Value = "8.870"
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("D27").Value = "0.1414"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").Value = "18.03"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("D30").Value = "1.401"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").Value = "1.475"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("D32").Value = "4.314"
$ws.Range("E32").Value = "  -4.23%  "
$ws.Range("D33").Value = "4.009"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "0.05145"
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("D35").Value = "1.167"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").Value = "0.7075"
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("D37").Value = "0.9973"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D39").Value = "0.01846"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").Value = "2.704"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").Value = "1.154.40"
$ws.Range("E41").Value = "  +4.91%  "
$ws.Range("D42").Value = "0.9221"
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("D43").Value = "5.957"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "0.4237"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("D45").Value = "70.04"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "103.04"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("D49").Value = "1.736"
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("D50").Value = "9.144"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "6.959"
$ws.Range("E51").Value = "  -0.55%  "
